# Update row 4 statistics values per May 2025 review 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 17.23014859704869
$ws.Range("F4").Value = 0.4327822234130386
$ws.Range("H4").Value = 16.87220457073091
$ws.Range("I4").Value = 17.24453245475632
$ws.Range("J4").Value = 17.6024764810741
$ws.Range("K4").Value = 17.60777899023485
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 19.43074970151516
$ws.Range("N4").Value = 2.327761773041207
$ws.Range("P4").Value = 18.14331576257711
$ws.Range("R4").Value = 19.67246997689321
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 17.52664365019696
$ws.Range("V4").Value = 2.244350442665643
$ws.Range("W4").Value = 15.78176450473163
$ws.Range("X4").Value = 15.97064659229483
$ws.Range("Y4").Value = 16.83168635074224
$ws.Range("Z4").Value = 18.38768340864436
$ws.Range("AB4").Value = 4
$ws.Range("AC4").Value = 0.2964950531482682
$ws.Range("AD4").Value = 1.903135609047713
$ws.Range("AE4").Value = -1.106591426760499
$ws.Range("AF4").Value = -0.869255256913644
$ws.Range("AG4").Value = -0.3840783885988142
$ws.Range("AH4").Value = 0.7816719214630981
$ws.Range("AJ4").Value = 4
$ws.Range("AK4").Value = 2.200601104466472
$ws.Range("AL4").Value = 2.610135201235659
$ws.Range("AN4").Value = 0.8523152630772808
$ws.Range("AP4").Value = 2.485254418480863
$ws.Range("AR4").Value = 4
$ws.Range("AS4").Value = 4.054655517009317
$ws.Range("AT4").Value = 1.903135609047713
$ws.Range("AU4").Value = 2.65156903710055
$ws.Range("AV4").Value = 2.888905206947405
$ws.Range("AW4").Value = 3.374082075262235
$ws.Range("AX4").Value = 4.539832385324147
$ws.Range("AZ4").Value = 4
$ws.Range("BA4").Value = 5.339835534230716
$ws.Range("BB4").Value = 2.610135201235659
$ws.Range("BD4").Value = 3.991549692841525
$ws.Range("BF4").Value = 5.624488848245107
$ws.Range("BH4").Value = 4
$ws.Range("BI4").Value = 0.09214871241436312
$ws.Range("BJ4").Value = 0.06601052259395453
$ws.Range("BL4").Value = 0.0568056970266501
$ws.Range("BM4").Value = 0.1002952374439998
$ws.Range("BN4").Value = 0.1356382528317129
$ws.Range("BO4").Value = 0.1591469002859662
$ws.Range("BP4").Value = 4
$ws.Range("BQ4").Value = 0.04717147691994795
$ws.Range("BR4").Value = 0.03402092871866552
$ws.Range("BT4").Value = 0.03531152072948598
$ws.Range("BU4").Value = 0.05188819912726571
$ws.Range("BV4").Value = 0.06374815531772768
$ws.Range("BX4").Value = 0.03300526129697726
$ws.Range("BY4").Value = 0.01701046435933276
